$d = $word.ActiveDocument

# --- Locate the end of the "Socks in the Dark" paragraph that ends with
# "It will be hard to get a specific color if you care which color you
# want for the day.  " (this paragraph currently also holds the
# "_GoBack" bookmark Word leaves at the last edited spot).
$marker = "It will be hard to get a specific color if you care which color you want for the day.  "

$findRng = $d.Content
$found = $findRng.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text for the Socks-in-the-Dark paragraph."
}

$endOfText = $findRng.End

# The "_GoBack" bookmark (auto-maintained by Word at the last edit point)
# currently sits right after this text, inside the same paragraph. The new
# paragraph being added should end up owning that bookmark instead, so drop
# it now and re-create it after the new paragraph is in place.
for ($i = $d.Bookmarks.Count; $i -ge 1; $i--) {
    $existing = $d.Bookmarks.Item($i)
    if ($existing.Name -eq "_GoBack") {
        $existing.Delete()
    }
}

# Split the paragraph right after the marker text, producing a new (empty)
# paragraph immediately following it.
$splitRng = $d.Range($endOfText, $endOfText)
$splitRng.InsertParagraphAfter()

# The freshly inserted paragraph mark sits one character after $endOfText.
$newParaMarkStart = $endOfText + 1

# Replace that paragraph mark's content/properties with the real new
# paragraph: List-Paragraph style, same numbered-list (numId 3) as the rest
# of the "Socks in the Dark" section, and the two runs of new text.
$newParaRng = $d.Range($newParaMarkStart, $newParaMarkStart + 1)
$newParaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>To solve the issue of the amount of socks in his drawer it would be a good idea to organize his choices into sections of his drawer so that he will always choose the right co</w:t></w:r><w:r><w:t>lor of sock with just one pull. Having no visual this person would also need to combine like colors of socks to avoid having that chance of picking a non matching pair of socks.</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@
$newParaRng.InsertXML($newParaXml)

# Re-create the "_GoBack" bookmark at the end of the text we just typed,
# matching Word's habit of tracking the most recent edit location.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bookmarkPos = $newPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

$d.Content.Find.Execute("", $false) | Out-Null
